$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3580
$ws.Range("I69").Value = 4090
$ws.Range("J69").Value = 2900
$ws.Range("K69").Value = 12270
$ws.Range("L69").Value = 8700
$ws.Range("M69").Value = -11396
$ws.Range("N69").Value = -10448

# Row 72
$ws.Range("H72").Value = 3580
$ws.Range("I72").Value = 4090
$ws.Range("J72").Value = 2900
$ws.Range("K72").Value = 36810
$ws.Range("L72").Value = 26100
$ws.Range("M72").Value = -32442
$ws.Range("N72").Value = -34836

# Row 80
$ws.Range("H80").Value = 1776.8077
$ws.Range("I80").Value = 654.875
$ws.Range("J80").Value = 3571.9
$ws.Range("K80").Value = 1964.625
$ws.Range("L80").Value = 10715.7
$ws.Range("M80").Value = -966.625
$ws.Range("N80").Value = -12711.7

# Row 83
$ws.Range("H83").Value = 1776.8077
$ws.Range("I83").Value = 654.875
$ws.Range("J83").Value = 3571.9
$ws.Range("K83").Value = 5893.875
$ws.Range("L83").Value = 32147.1
$ws.Range("M83").Value = -901.875
$ws.Range("N83").Value = -42131.10000000001

# Row 107
$ws.Range("H107").Value = 3709.842
$ws.Range("I107").Value = 3882.2778
$ws.Range("J107").Value = 606
$ws.Range("K107").Value = 3882.2778
$ws.Range("L107").Value = 606
$ws.Range("M107").Value = -1962.2778
$ws.Range("N107").Value = -4446

# Row 113
$ws.Range("H113").Value = 2453.7144
$ws.Range("I113").Value = 2420
$ws.Range("J113").Value = 2498.6667
$ws.Range("K113").Value = 2420
$ws.Range("L113").Value = 2498.6667
$ws.Range("M113").Value = 834
$ws.Range("N113").Value = -9006.6667

# Row 132
$ws.Range("H132").Value = 3967.5
$ws.Range("I132").Value = 4610.516
$ws.Range("J132").Value = 2155.3635
$ws.Range("K132").Value = 13831.548
$ws.Range("L132").Value = 6466.0905
$ws.Range("M132").Value = -11301.548
$ws.Range("N132").Value = -11526.0905

# Row 137
$ws.Range("H137").Value = 1442.2307
$ws.Range("I137").Value = 615.3333
$ws.Range("J137").Value = 1690.3
$ws.Range("K137").Value = 1845.9999
$ws.Range("L137").Value = 5070.9
$ws.Range("M137").Value = 704.0001
$ws.Range("N137").Value = -10170.9

# Row 141
$ws.Range("H141").Value = 3047.4666
$ws.Range("I141").Value = 1529.9584
$ws.Range("J141").Value = 9117.5
$ws.Range("K141").Value = 4589.8752
$ws.Range("L141").Value = 27352.5
$ws.Range("M141").Value = 590.1247999999996

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 407317.38
$ws.Range("I32").Value = 521469.22
$ws.Range("J32").Value = 14920.375
$ws.Range("K32").Value = 521469.22
$ws.Range("L32").Value = 14920.375
$ws.Range("M32").Value = -521182.22
$ws.Range("N32").Value = -15494.375

# Row 61
$ws.Range("H61").Value = 14495535
$ws.Range("I61").Value = 25643328
$ws.Range("J61").Value = 3406
$ws.Range("K61").Value = 25643328
$ws.Range("L61").Value = 3406
$ws.Range("M61").Value = -25643116
$ws.Range("N61").Value = -3830

# Row 74
$ws.Range("H74").Value = 980.4091
$ws.Range("I74").Value = 875.75
$ws.Range("J74").Value = 1106
$ws.Range("K74").Value = 875.75
$ws.Range("L74").Value = 1106
$ws.Range("M74").Value = -1.75
$ws.Range("N74").Value = -2854

# Row 77
$ws.Range("H77").Value = 980.4091
$ws.Range("I77").Value = 875.75
$ws.Range("J77").Value = 1106
$ws.Range("K77").Value = 4378.75
$ws.Range("L77").Value = 5530
$ws.Range("M77").Value = -10.75
$ws.Range("N77").Value = -14266

# Row 110
$ws.Range("H110").Value = 5370.3335
$ws.Range("I110").Value = 5190.4287
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 5190.4287
$ws.Range("L110").Value = 6000
$ws.Range("M110").Value = -3145.4287
$ws.Range("N110").Value = -10090

# Row 132
$ws.Range("H132").Value = 6038.8667
$ws.Range("I132").Value = 7548.4287
$ws.Range("J132").Value = 4718
$ws.Range("K132").Value = 22645.2861
$ws.Range("L132").Value = 14154
$ws.Range("M132").Value = -20115.2861
$ws.Range("N132").Value = -19214

# Row 136
$ws.Range("H136").Value = 14495535
$ws.Range("I136").Value = 25643328
$ws.Range("J136").Value = 3406
$ws.Range("K136").Value = 76929984
$ws.Range("L136").Value = 10218
$ws.Range("M136").Value = -76927434
$ws.Range("N136").Value = -15318

$ws = $wb.Worksheets.Item("BSM")
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Row 134
$ws.Range("H134").Value = 2201.3333
$ws.Range("I134").Value = 1927.3549
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 5782.0647
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -3247.0647
$ws.Range("N134").Value = -16770

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1574

# Row 31
$ws.Range("H31").Value = 1352.0944
$ws.Range("I31").Value = 1071.1621
$ws.Range("J31").Value = 2001.75
$ws.Range("K31").Value = 1071.1621
$ws.Range("L31").Value = 2001.75
$ws.Range("M31").Value = -776.1621
$ws.Range("N31").Value = -2591.75

# Row 34
$ws.Range("H34").Value = 1352.0944
$ws.Range("I34").Value = 1071.1621
$ws.Range("J34").Value = 2001.75
$ws.Range("K34").Value = 1071.1621
$ws.Range("L34").Value = 2001.75
$ws.Range("M34").Value = -869.1621
$ws.Range("N34").Value = -2405.75

# Row 58
$ws.Range("H58").Value = 2869.6667
$ws.Range("I58").Value = 2902.1428
$ws.Range("J58").Value = 2824.2
$ws.Range("K58").Value = 2902.1428
$ws.Range("L58").Value = 2824.2
$ws.Range("M58").Value = -2699.1428
$ws.Range("N58").Value = -3230.2

# Row 99
$ws.Range("H99").Value = 1600
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = 498

# Row 113
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5340

# Row 126
$ws.Range("H126").Value = 1600
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -530

# Row 136
$ws.Range("H136").Value = 2869.6667
$ws.Range("I136").Value = 2902.1428
$ws.Range("J136").Value = 2824.2
$ws.Range("K136").Value = 8706.428400000001
$ws.Range("L136").Value = 8472.599999999999
$ws.Range("M136").Value = -6156.428400000001
$ws.Range("N136").Value = -13572.6

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 901.0941
$ws.Range("I68").Value = 1296.6666
$ws.Range("J68").Value = 886.62195
$ws.Range("K68").Value = 3889.9998
$ws.Range("L68").Value = 2659.86585
$ws.Range("M68").Value = -3078.9998
$ws.Range("N68").Value = -4281.86585

# Row 71
$ws.Range("H71").Value = 901.0941
$ws.Range("I71").Value = 1296.6666
$ws.Range("J71").Value = 886.62195
$ws.Range("K71").Value = 11669.9994
$ws.Range("L71").Value = 7979.59755
$ws.Range("M71").Value = -7613.999400000001
$ws.Range("N71").Value = -16091.59755

# Row 82
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30812
$ws.Range("M82").ClearContents()

# Row 85
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32808
$ws.Range("M85").ClearContents()

# Row 107
$ws.Range("H107").Value = 2121.8108
$ws.Range("I107").Value = 244.42857
$ws.Range("J107").Value = 2559.8667
$ws.Range("K107").Value = 733.28571
$ws.Range("L107").Value = 7679.6001
$ws.Range("M107").Value = 1186.71429
$ws.Range("N107").Value = -11519.6001

# Row 137
$ws.Range("H137").Value = 9999.066000000001
$ws.Range("I137").Value = 13128.6
$ws.Range("J137").Value = 3740
$ws.Range("K137").Value = 39385.8
$ws.Range("L137").Value = 11220
$ws.Range("M137").Value = -34285.8
$ws.Range("N137").Value = -21420

$ws = $wb.Worksheets.Item("GSM")
# Row 42
$ws.Range("H42").Value = 55430
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 55430
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 55430
$ws.Range("N42").Value = -56400

# Row 107
$ws.Range("H107").Value = 609.75
$ws.Range("I107").Value = 428.66666
$ws.Range("J107").Value = 842.5714
$ws.Range("K107").Value = 428.66666
$ws.Range("L107").Value = 842.5714
$ws.Range("M107").Value = 1491.33334
$ws.Range("N107").Value = -4682.5714

# Row 115
$ws.Range("H115").Value = 55430
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 55430
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 55430
$ws.Range("N115").Value = -57780

# Row 132
$ws.Range("H132").Value = 3241.9092
$ws.Range("I132").Value = 3048
$ws.Range("J132").Value = 3522
$ws.Range("K132").Value = 9144
$ws.Range("L132").Value = 10566
$ws.Range("M132").Value = -6614

# Row 136
$ws.Range("H136").Value = 79326
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 79326
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 237978
$ws.Range("N136").Value = -243078

$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 40097.25
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40097.25
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40097.25
$ws.Range("N92").Value = -45089.25

$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 50121.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 50121.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 50121.5
$ws.Range("N82").Value = -50887.5

# Row 85
$ws.Range("H85").Value = 50121.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 50121.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 50121.5
$ws.Range("N85").Value = -52773.5

# Row 137
$ws.Range("H137").Value = 67857.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 67857.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 67857.5
$ws.Range("N137").Value = -78057.5
